$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "69.641.70"
$ws.Range("E2").Value = "  -0.50%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.811.83"
$ws.Range("E3").Value = "  +1.97%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'611.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.91%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'176.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.06%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.808.61"
$ws.Range("E7").Value = "  +1.95%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.00%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -1.80%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.49%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "'6.47"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.53%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "'0.481"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.87%  "

# Row 13 - Avalanche
$ws.Range("D13").Value = "'39.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.84%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  -2.31%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.447.98"
$ws.Range("E15").Value = "  +1.97%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.815.82"
$ws.Range("E16").Value = "  +2.16%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "69.678.79"
$ws.Range("E17").Value = "  -0.43%  "

# Row 18 - Polkadot
$ws.Range("E18").Value = "  -0.48%  "

# Row 19 - TRON
$ws.Range("E19").Value = "  -3.27%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "'16.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.18%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'506.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.13%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'9.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.59%  "

# Row 23 - Polygon
$ws.Range("E23").Value = "  +2.16%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'86.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.67%  "

# Row 25 - Fetch.AI
$ws.Range("E25").Value = "  -3.28%  "

# Row 26 - PEPE
$ws.Range("D26").Value = "'0.0000142"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.23%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("E27").Value = "  -3.76%  "

# Row 28 - RenderToken
$ws.Range("E28").Value = "  -6.83%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  +0.05%  "

# Row 30 - ImmutableX
$ws.Range("E30").Value = "  +1.78%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +1.42%  "

# Row 32 - NEARProtocol
$ws.Range("E32").Value = "  +0.31%  "

# Row 33 - EthereumClassic
$ws.Range("D33").Value = "'31.68"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.67%  "

# Row 34 - Hedera
$ws.Range("E34").Value = "  -1.65%  "

# Row 35 - FirstDigitalUSD
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.10%  "

# Row 36 - Mantle
$ws.Range("D36").Value = "'1.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.42%  "

# Row 37 - Filecoin
$ws.Range("E37").Value = "  -1.84%  "

# Row 38 - Kaspa
$ws.Range("E38").Value = "  +5.16%  "

# Row 39 - Bittensor
$ws.Range("D39").Value = "'482.61"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +13.49%  "

# Row 40 - TheGraph
$ws.Range("E40").Value = "  -0.12%  "

# Row 41 - dogwifhat
$ws.Range("E41").Value = "  +6.76%  "

# Row 42 - Stacks
$ws.Range("E42").Value = "  -3.11%  "

# Row 43 - OKB
$ws.Range("D43").Value = "'49.72"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.46%  "

# Row 44 - Arweave
$ws.Range("D44").Value = "'43.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.20%  "

# Row 45 - Cosmos
$ws.Range("E45").Value = "  -1.98%  "

# Row 46 - Maker
$ws.Range("D46").Value = "2.917.08"
$ws.Range("E46").Value = "  -2.71%  "

# Row 47 - VeChain
$ws.Range("E47").Value = "  -1.06%  "

# Row 48 - Monero
$ws.Range("D48").Value = "'139.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.21%  "

# Row 49 - USDe
$ws.Range("E49").Value = "  +0.04%  "

# Row 50 - InjectiveProtocol
$ws.Range("D50").Value = "'26.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.47%  "
